$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "Get the shipping boxes..." paragraph gains a second sentence as a
# brand-new run (not merged into the existing run).
# ---------------------------------------------------------------------------
$marker1 = "The boxes are labelled by the instrument they are for."
$searchRange1 = $d.Content
$searchRange1.Find.Execute($marker1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterMarker1 = $searchRange1.End

$newText1 = " If there is not a labelled box for your instrument and no other box in the hot lab will fit your instrument, you may get a box from the recycling pile outside the infusion room. Check in the afternoon for the best available selection."

$insertRange1 = $d.Range($afterMarker1, $afterMarker1)
$insertRange1.InsertAfter($newText1)

# Force the newly inserted text into its own run (rather than being merged
# into the preceding run) by briefly dropping a bookmark at the seam and
# removing it again -- this splits the underlying run without leaving any
# bookmark behind.
$seamRange1 = $d.Range($afterMarker1, $afterMarker1)
$d.Bookmarks.Add("ZZZ_TempSplit1", $seamRange1) | Out-Null
$d.Bookmarks("ZZZ_TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Part 2: "Place each instrument..." paragraph gets a new lead-in, split
# across three runs with the (moved) _GoBack bookmark sitting between the
# first and second new runs.
# ---------------------------------------------------------------------------
$searchRange2 = $d.Content
$searchRange2.Find.Execute("Place each instrument", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pStart = $searchRange2.Start

$runA = "If your box has no shaped Styrofoam, pack the instrument with bubble wrap and/or Styrofoam; K&S will send the box back with sh"
$runB = "aped Styrofoam. Otherwise, p"

# Replace the leading capital "P" of "Place" with runA + runB (ending in the
# lower-case "p" that used to be the "P").
$pCharRange = $d.Range($pStart, $pStart + 1)
$pCharRange.Text = $runA + $runB

$splitPos1 = $pStart + $runA.Length
$splitPos2 = $pStart + $runA.Length + $runB.Length

# Split off run3 ("lace each instrument...") first, at the higher offset, so
# the earlier offset ($splitPos1) stays valid afterwards.
$seamRange2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("ZZZ_TempSplit2", $seamRange2) | Out-Null
$d.Bookmarks("ZZZ_TempSplit2").Delete()

# ---------------------------------------------------------------------------
# Part 3: the stray _GoBack bookmark that used to sit at the very end of the
# document is removed, and a fresh one is inserted between run1/run2 of part
# 2 above (i.e. it effectively moves).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$seamRange1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("_GoBack", $seamRange1) | Out-Null
